$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set cell values in the order the original strings were authored so that
# the shared-string table indices line up with the target workbook:
# Issue(30), Not able... (31), Resolving method(32), The pretrained...(33)
$ws.Range("E1").Value = "Issue"
$ws.Range("E2").Value = "Not able to download the weights using the command used in Keras website"
$ws.Range("F1").Value = "Resolving method"
$ws.Range("F2").Value = "The pretrained model weights are downloaded to the kaggle kernel"

# The OOXML <col> width is derived from ColumnWidth through the host's
# character->pixel conversion (quantized to 2 decimal places of
# ColumnWidth, then to whole pixels). These inputs are the closest
# achievable match to the target stored widths (70.42578125 / 17.28515625).
$ws.Columns.Item(5).ColumnWidth = 69.59
$ws.Columns.Item(6).ColumnWidth = 16.42

$ws.Range("F2").Select()
